# Updates cryptos list data: Price (D) and Volume(1h) (E) columns,
# plus a couple of row swaps/renames in Coin (B) / Link (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as plain text (e.g. "70.728.87",
# "492.30"). Excel auto-detects numeric-looking text on assignment
# to Value, which would silently drop significant trailing zeros
# (e.g. "492.30" -> 492.3). Force text storage via NumberFormat "@"
# around the write, then ClearFormats() so no stray cell style is
# left behind (the source cells carry no explicit style).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) "70.728.87"
$ws.Cells.Item(2, 5).Value = "  +7.36%  "

Set-TextValue $ws.Cells.Item(3, 4) "3.633.51"
$ws.Cells.Item(3, 5).Value = "  +7.17%  "

$ws.Cells.Item(4, 5).Value = "  +0.00%  "

Set-TextValue $ws.Cells.Item(5, 4) "593.04"
$ws.Cells.Item(5, 5).Value = "  +5.00%  "

Set-TextValue $ws.Cells.Item(6, 4) "191.43"
$ws.Cells.Item(6, 5).Value = "  +8.34%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.647"
$ws.Cells.Item(7, 5).Value = "  +2.59%  "

Set-TextValue $ws.Cells.Item(8, 4) "3.617.92"
$ws.Cells.Item(8, 5).Value = "  +6.87%  "

$ws.Cells.Item(9, 5).Value = "  -0.05%  "

$ws.Cells.Item(10, 5).Value = "  +2.80%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.662"
$ws.Cells.Item(11, 5).Value = "  +4.32%  "

Set-TextValue $ws.Cells.Item(12, 4) "57.87"
$ws.Cells.Item(12, 5).Value = "  +7.35%  "

Set-TextValue $ws.Cells.Item(13, 4) "0.0000294"
$ws.Cells.Item(13, 5).Value = "  +5.45%  "

Set-TextValue $ws.Cells.Item(14, 4) "9.74"
$ws.Cells.Item(14, 5).Value = "  +5.37%  "

Set-TextValue $ws.Cells.Item(15, 4) "4.211.22"
$ws.Cells.Item(15, 5).Value = "  +7.22%  "

Set-TextValue $ws.Cells.Item(16, 4) "3.626.04"
$ws.Cells.Item(16, 5).Value = "  +7.19%  "

Set-TextValue $ws.Cells.Item(17, 4) "19.38"
$ws.Cells.Item(17, 5).Value = "  +6.28%  "

Set-TextValue $ws.Cells.Item(18, 4) "70.650.33"
$ws.Cells.Item(18, 5).Value = "  +7.38%  "

Set-TextValue $ws.Cells.Item(19, 4) "12.59"
$ws.Cells.Item(19, 5).Value = "  +5.79%  "

$ws.Cells.Item(20, 5).Value = "  +0.87%  "

$ws.Cells.Item(21, 5).Value = "  +5.25%  "

Set-TextValue $ws.Cells.Item(22, 4) "492.30"
$ws.Cells.Item(22, 5).Value = "  +6.00%  "

Set-TextValue $ws.Cells.Item(23, 4) "5.50"
$ws.Cells.Item(23, 5).Value = "  +11.18%  "

Set-TextValue $ws.Cells.Item(24, 4) "16.80"
$ws.Cells.Item(24, 5).Value = "  +14.88%  "

$ws.Cells.Item(25, 5).Value = "  +9.02%  "

Set-TextValue $ws.Cells.Item(26, 4) "90.84"
$ws.Cells.Item(26, 5).Value = "  +1.42%  "

$ws.Cells.Item(27, 5).Value = "  +5.73%  "

Set-TextValue $ws.Cells.Item(28, 4) "11.22"
$ws.Cells.Item(28, 5).Value = "  +5.30%  "

Set-TextValue $ws.Cells.Item(29, 4) "9.38"
$ws.Cells.Item(29, 5).Value = "  +7.46%  "

Set-TextValue $ws.Cells.Item(30, 4) "32.23"
$ws.Cells.Item(30, 5).Value = "  +3.45%  "

Set-TextValue $ws.Cells.Item(31, 4) "7.67"
$ws.Cells.Item(31, 5).Value = "  +15.96%  "

Set-TextValue $ws.Cells.Item(32, 4) "12.25"
$ws.Cells.Item(32, 5).Value = "  +6.66%  "

Set-TextValue $ws.Cells.Item(33, 4) "618.72"
$ws.Cells.Item(33, 5).Value = "  +6.47%  "

$ws.Cells.Item(34, 5).Value = "  +8.00%  "

Set-TextValue $ws.Cells.Item(35, 4) "65.30"
$ws.Cells.Item(35, 5).Value = "  +4.73%  "

Set-TextValue $ws.Cells.Item(36, 4) "0.0₃0831"
$ws.Cells.Item(36, 5).Value = "  +10.98%  "

$ws.Cells.Item(37, 5).Value = "  +3.93%  "

$ws.Cells.Item(38, 2).Value = "TheGraph"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Cells.Item(38, 4) "0.404"
$ws.Cells.Item(38, 5).Value = "  +7.00%  "

$ws.Cells.Item(39, 2).Value = "Dai"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Cells.Item(39, 4) "1.00"
$ws.Cells.Item(39, 5).Value = "  -0.08%  "

Set-TextValue $ws.Cells.Item(40, 4) "37.93"
$ws.Cells.Item(40, 5).Value = "  +5.20%  "

Set-TextValue $ws.Cells.Item(41, 4) "3.64"
$ws.Cells.Item(41, 5).Value = "  +1.34%  "

Set-TextValue $ws.Cells.Item(42, 4) "3.352.89"
$ws.Cells.Item(42, 5).Value = "  +7.92%  "

Set-TextValue $ws.Cells.Item(43, 4) "3.08"
$ws.Cells.Item(43, 5).Value = "  +8.22%  "

Set-TextValue $ws.Cells.Item(44, 4) "0.0448"
$ws.Cells.Item(44, 5).Value = "  +7.30%  "

$ws.Cells.Item(45, 5).Value = "  +9.02%  "

Set-TextValue $ws.Cells.Item(46, 4) "3.38"
$ws.Cells.Item(46, 5).Value = "  +5.97%  "

Set-TextValue $ws.Cells.Item(47, 4) "0.138"
$ws.Cells.Item(47, 5).Value = "  +2.59%  "

Set-TextValue $ws.Cells.Item(48, 4) "9.23"
$ws.Cells.Item(48, 5).Value = "  +8.42%  "

Set-TextValue $ws.Cells.Item(49, 4) "2.74"
$ws.Cells.Item(49, 5).Value = "  +6.50%  "

Set-TextValue $ws.Cells.Item(50, 4) "3.36"
$ws.Cells.Item(50, 5).Value = "  +6.25%  "

$ws.Cells.Item(51, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Cells.Item(51, 4) "0.998"
$ws.Cells.Item(51, 5).Value = "  -0.06%  "
